$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Formula = "=-27 -79.29 -27 -79.29 -79.29 -79.29"
$ws.Range("F3").Formula = "=-40.94 -79.29 -79.29 -40.94 -79.29 -40.94 -79.29 -40.94 -74.07 -79.29 -27 -40.94 -27 -79.29 -40.94 -27.00 -74.07"
$ws.Range("H3").Formula = "=-79.29 -74.07 -40.94 -27 -27 -27 -74.07 -79.29 -74.07 -40.94 -74.07 -40.94"

# Row 4
$ws.Range("E4").Formula = "=-67.84 -297.51 -67.84 -67.84 -67.84 -21.89 -39.84 -67.84 -21.89 -67.84 -21.89 -21.89 -39.84 -67.84"
$ws.Range("F4").Formula = "=-39.84 -21.89 -297.51 -67.84 -21.89 -39.84 -297.51 -39.84 -297.51 -297.51 -21.89"
$ws.Range("H4").Formula = "=-67.84 -297.51 -297.51 -67.84 -67.84 -297.51"
$ws.Range("I4").Formula = "=-21.89 -297.51 -39.84"

$excel.Calculate()
